$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 1237.087816454243
$ws.Range("D2").Value = 2172.779588794647
